# issue #5: stock data from json to db
# The "股票" (stock) worksheet gains three new columns describing where the
# record came from: category (between property_category and date),
# source_file and index (appended after legislator_id).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (I) so the layout
# becomes: ... property_category, category, date, legislator_name,
# legislator_id, ...  The insert also carries the header/data formatting
# from the column immediately to its left.
$ws.Columns("I").Insert()

# -- Header row --------------------------------------------------------
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# -- Data rows -----------------------------------------------------------
$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"
$ws.Range("I4").Value = "normal"

$ws.Range("M2").Value = "tmp80511"
$ws.Range("M3").Value = "tmp80511"
$ws.Range("M4").Value = "tmp80511"

$ws.Range("N2").Value = 90
$ws.Range("N3").Value = 91
$ws.Range("N4").Value = 92

# Match the formatting already used for the rest of the header/data columns
# by copying it from the neighbouring column instead of re-building it
# property by property (avoids creating redundant style entries).
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)

$ws.Range("L2:L4").Copy()
$ws.Range("M2:N4").PasteSpecial(-4122)
